$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B; the old column B shifts to D
# and the old column C shifts to E.
$ws.Range("B:C").Insert()

# New header labels for the two freshly inserted columns.
$ws.Range("B1").Value2 = "Jun_17"
$ws.Range("C1").Value2 = "Jun_15"

# The new columns get the same "UN" placeholder text used throughout column B/C.
$ws.Range("B2:C27").Value2 = "UN"

# Match column widths (C, D, E all keep the same custom width as the
# original column C did).
$ws.Columns.Item(3).ColumnWidth = 7.14
$ws.Columns.Item(4).ColumnWidth = 7.14
$ws.Columns.Item(5).ColumnWidth = 7.14
